$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column C for rows 2-8 from
# serial 45207 (2023-10-08) to serial 45208 (2023-10-09).
for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 3).Value = 45208
}
